$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the two obsolete single-purpose UOM test sheets.
$wb.Worksheets.Item("tc_UOM_004").Delete()
$wb.Worksheets.Item("tc_UOM_009").Delete()

# Complete the tc_UOM_005_010 subset into tc_UOM_004_005_010:
# insert a new "EditFieldOfUOM" column before the existing "Updated UOM name"
# column and refresh a couple of automation data values.
$ws = $wb.Worksheets.Item("tc_UOM_005_010")
$ws.Columns.Item(5).Insert()
$ws.Range("E1").Value = "EditFieldOfUOM"
$ws.Range("E2").Value = "*UOM Name,UOM Descriptions"
$ws.Range("B2").Value = "Pauto"
$ws.Range("F2").Value = "PAutoTest2"
$ws.Name = "tc_UOM_004_005_010"

$ws.Activate()
$ws.Range("F2").Select()
